# Berichten.xlsx edit:
#  - Remove the two "EasyRoads" notification rows (Nummer.../Smash?,
#    FlobberBark/Fourtnite, Snoepjesgever/Iemand een snoepje?)
#  - Add a batch of extra "irrelevant" social-media style notifications
#  - Table grows from 17 data rows (A1:E18) to 25 data rows (A1:E26)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Grow the XML-mapped table + its autofilter to the new extent first so
# the newly written rows inherit the table formatting/range.
$lo.Resize($ws.Range("A1:E26"))

# r -> (Name, Message, POI, Platform, Image)
$rows = @{
    2  = @("Richard", "Is dit hem?", 0, "Twitter", "IsDitHem")
    3  = @("Dennis", "Ik heb gister een broodje pindakaas gegeten. #peanutbutter4lyfe #lekker", 0, "Twitter", "")
    4  = @("Nummer 06000000001", "Waar kan ik de aphotheek vinden?", 0, "Twitter", "")
    5  = @("LisaNL12", "Lekker dagje winkelen met mij vriendinnen!!", 0, "Twitter", "")
    6  = @("Jordaan038", "kopje koffie? Ik krijg er nooit genoeg van!", 0, "Twitter", "")
    7  = @("Fluffy06", "Zag net zo'n schattig knuffeltje. Heb er één gekocht voor mijn broertje maar wil er eigelijk ook zelf één XD", 0, "Twitter", "")
    8  = @("Roos", "Eindelijk klaar met leren! Nu tijd om te relaxen! ", 0, "Whatsapp", "")
    9  = @("Kees1999", "Dat moment wanneer je niet weet wat je moet tweeten #random", 0, "Twitter", "")
    10 = @("FlobberBark", "Wat een mooie dag is het vandaag! Lekker een luchte scheppen :D", 0, "Facebook", "")
    11 = @("Natasha", "Whoa kan niet geloven hoe goedkoop die schoenen waren! #lucky #opruiming", 0, "Whatsapp", "")
    12 = @("Ricky4", "Vandaag ga ik lekker niks doen… behalve dit bericht poste haha", 0, "Twitter", "")
    13 = @("Marloesje", "Naar de kapper staat op de planning vandaag! Is ook wel nodig haha", 0, "Twitter", "")
    14 = @("Lisa", "Er is hier net een overval gepleegd!", 1, "SMS", "")
    15 = @("Lars", "Ik zie iemand die zich verdacht gedraagt.", 2, "Facebook", "Walking")
    16 = @("Niels", "Ik zag hier net een persoon voorbij lopen die er uitziet als die overvaller.", 2, "Whatsapp", "")
    17 = @("GentStudent", "#DoeHetVoorGentStudent42.2", 0, "Twitter", "")
    18 = @("VervelendeZeurpiet", "De bus is weer laat @Synthus #sad #late", 0, "Twitter", "")
    19 = @("Tom", "ik heb net een straatroof gezien, er was een vrouw overvallen door een man.", 1, "Bellen", "")
    20 = @("Geert", "Heb iemand gezien met die verdachteomschrijving.", 3, "Bellen", "")
    21 = @("Sweet Tooth", "Zo hyped!! Buitelands snoep komt morgen binnen!! ", 0, "Facebook", "")
    22 = @("Sandra", "Wow, straatroof in het daglicht! Waar gaat de wereld naartoe? #sad", 1, "Twitter", "Running")
    23 = @("Femke", "Hier staat iemand die veel op het signalement lijkt al een tijdje te wachten.", 3, "Whatsapp", "Suspect")
    24 = @("Anna", "Hier liep wel iemand die erop leek.", 3, "Twitter", "")
    25 = @("Tim", "Uit mijn raam zag ik de verdachte.", 2, "Bellen", "")
    26 = @("Iris", "Ik zie hem hier.", 1, "Twitter", "")
}

foreach ($r in 2..26) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
    if ($vals[4] -ne "") {
        $ws.Cells.Item($r, 5).Value = $vals[4]
    } else {
        $ws.Cells.Item($r, 5).Value = ""
    }
}

# Make sure every data row carries the same text-formatted style the
# original rows used for columns A, B, D, E (column C already inherits
# the sheet's style-1 text format from the column definition).
$ws.Range("A2:B26").NumberFormat = "@"
$ws.Range("D2:E26").NumberFormat = "@"

# Restore the cursor position recorded in the saved file.
$ws.Range("D12").Select()
